# Modify cases for no survey selectivity if no survey data
# Update the "E" scenario-code column (survey selectivity) for the rows that
# represent "no survey data" cases, replacing the old shared codes
# (E0, E1, E2, E3) with new dedicated codes (E100, E101, E102, E103) so that
# they no longer alias onto the regular E0-E3 survey-selectivity cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios")

# Rows previously coded E0 (string index 23) -> E100
$ws.Range("E2").Value  = "E100"
$ws.Range("E4").Value  = "E100"
$ws.Range("E26").Value = "E100"
$ws.Range("E28").Value = "E100"

# Rows previously coded E1 (string index 19) -> E101
$ws.Range("E9").Value  = "E101"
$ws.Range("E11").Value = "E101"
$ws.Range("E33").Value = "E101"
$ws.Range("E35").Value = "E101"

# Rows previously coded E2 (string index 3) -> E102
$ws.Range("E16").Value = "E102"
$ws.Range("E50").Value = "E102"

# Rows previously coded E3 (string index 63) -> E103
$ws.Range("E21").Value = "E103"
$ws.Range("E65").Value = "E103"

# Update the active-cell selection recorded for the scenarios sheet view.
$ws.Activate()
$ws.Range("H10").Select()

# Re-apply the AutoFilter (toggle off/on) so that the stale sortState that
# had been captured inside the autoFilter element is cleared, while leaving
# the autoFilter itself (and the separate standalone sortState below it)
# intact.
$ws.Range("A1:N69").AutoFilter() | Out-Null
$ws.Range("A1:N69").AutoFilter() | Out-Null

$wb.Save()
